$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("地方台JS脚本")

# Insert a new row at position 43 (shifts existing rows 43-61 down to 44-62,
# and auto-extends the merged cells / dimension that span across the break).
$ws.Rows.Item(43).Insert()

# Copy the formatting (borders/alignment/font) from row 42 (the "shanxi.js"
# row) onto the freshly inserted row 43 so the new row matches the rest of
# the merged "山西" (Shanxi) block.
$ws.Range("A42:D42").Copy()
$ws.Range("A43:D43").PasteSpecial(-4122)

# Fill in the new script entry.
$ws.Range("B43").Value = "shanxi_new.js"
$ws.Range("C43").Value = "本地"

# D41:D42 used to hold the "卫视,山西省频道" note for shenzhen.js/shanxi.js;
# extend that merged note down through the new row D43 as well.
$ws.Range("D41:D43").Merge()
$ws.Range("D41:D43").HorizontalAlignment = -4131

# The autofilter needs to grow from D61 to D62 - toggle it off then back on
# over the new full range (re-toggling an already-on filter would just turn
# it off).
$ws.AutoFilterMode = $false
$ws.Range("A1:D62").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name in sync with the new
# autofilter range.
$fd = $ws.Names.Item("_xlnm._FilterDatabase")
$fd.RefersTo = "=地方台JS脚本!`$A`$1:`$D`$62"
